{"js": "// The edit:\n//  1. Removes the \"_GoBack\" bookmark that currently sits at the end of the\n//     JWT_SECRET paragraph.\n//  2. Merges the three paragraphs that hold the (word-wrapped) OPENAI_API_KEY\n//     value into a single paragraph, re-assembling the full key text.\n//  3. Re-inserts the \"_GoBack\" bookmark inside that merged paragraph, right\n//     before the final \"T2yXKwcA\" fragment of the key.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraphs we need by their (unique) text content instead of\n// hard-coded indexes, so the script is resilient to minor positional shifts.\nlet openAiParagraph = null;\nlet middleParagraph = null;\nlet tailParagraph = null;\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"OPENAI_API_KEY=\") === 0) {\n    openAiParagraph = items[i];\n    middleParagraph = items[i + 1];\n    tailParagraph = items[i + 2];\n    break;\n  }\n}\n\nif (!openAiParagraph || !middleParagraph || !tailParagraph) {\n  throw new Error(\"Could not locate the OPENAI_API_KEY paragraphs.\");\n}\n\nconst middleText = middleParagraph.text; // e.g. \"wpPAA19...qtM3u\"\nconst tailText = tailParagraph.text; // e.g. \"wnA9oOT2yXKwcA\"\n\n// The tail paragraph's text is itself the concatenation of the last OpenAI\n// key fragment (\"wnA9oO\") and the short trailing fragment that stays after\n// the relocated bookmark (\"T2yXKwcA\"). Split it accordingly.\nconst bookmarkSuffix = \"T2yXKwcA\";\nconst tailKeyFragment = tailText.slice(0, tailText.length - bookmarkSuffix.length);\n\n// 1. Remove the bookmark from its current location (end of JWT_SECRET line).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Append the middle + tail key fragments onto the OPENAI_API_KEY paragraph.\nopenAiParagraph.getRange(\"End\").insertText(middleText + tailKeyFragment, Word.InsertLocation.end);\n\n// 3. Drop the now-redundant paragraphs (their text has been folded in above).\nmiddleParagraph.delete();\ntailParagraph.delete();\nawait context.sync();\n\n// 4. Re-insert the bookmark right after the key, then append the final\n//    trailing fragment after it.\nopenAiParagraph.getRange(\"End\").insertBookmark(\"_GoBack\");\nopenAiParagraph.getRange(\"End\").insertText(bookmarkSuffix, Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# The edit:\n#  1. Removes the \"_GoBack\" bookmark that currently sits at the end of the\n#     JWT_SECRET paragraph.\n#  2. Merges the three paragraphs that hold the (word-wrapped) OPENAI_API_KEY\n#     value into a single paragraph, re-assembling the full key text.\n#  3. Re-inserts the \"_GoBack\" bookmark inside that merged paragraph, right\n#     before the final \"T2yXKwcA\" fragment of the key.\n\n$d = $word.ActiveDocument\n\n# --- locate the three paragraphs by content, not by hard-coded index -------\n$openAiIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"OPENAI_API_KEY=\")) {\n        $openAiIndex = $i\n        break\n    }\n}\nif ($openAiIndex -eq -1) {\n    throw \"Could not locate the OPENAI_API_KEY paragraph.\"\n}\n\n$middleParagraph = $d.Paragraphs.Item($openAiIndex + 1)\n$tailParagraph = $d.Paragraphs.Item($openAiIndex + 2)\n\n$middleText = $middleParagraph.Range.Text.TrimEnd([char]13, [char]7)\n$tailText = $tailParagraph.Range.Text.TrimEnd([char]13, [char]7)\n\n# The tail paragraph's text is itself the concatenation of the last OpenAI\n# key fragment (\"wnA9oO\") and the short trailing fragment that stays after\n# the relocated bookmark (\"T2yXKwcA\"). Split it accordingly.\n$bookmarkSuffix = \"T2yXKwcA\"\n$tailKeyFragment = $tailText.Substring(0, $tailText.Length - $bookmarkSuffix.Length)\n\n# 1. Remove the bookmark from its current location (end of JWT_SECRET line).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Append the middle + tail key fragments onto the OPENAI_API_KEY paragraph.\n$openAiParagraph = $d.Paragraphs.Item($openAiIndex)\n$r = $openAiParagraph.Range\n[void]$r.MoveEnd(1, -1)   # exclude the trailing paragraph mark\n$r.Collapse(0)            # 0 = wdCollapseEnd\n$r.InsertAfter($middleText + $tailKeyFragment)\n\n# 3. Drop the now-redundant paragraphs (their text has been folded in above).\n#    After each delete, the following paragraph shifts up to the same index.\n$d.Paragraphs.Item($openAiIndex + 1).Range.Delete()\n$d.Paragraphs.Item($openAiIndex + 1).Range.Delete()\n\n# 4. Append the final trailing fragment after the key.\n$openAiParagraph2 = $d.Paragraphs.Item($openAiIndex)\n$r2 = $openAiParagraph2.Range\n[void]$r2.MoveEnd(1, -1)\n$r2.Collapse(0)\n$r2.InsertAfter($bookmarkSuffix)\n\n# 5. Re-insert the bookmark right before that trailing fragment. (Collapsing\n#    exactly on the paragraph-mark boundary is avoided by moving back from a\n#    safely-interior point instead of collapsing the full-paragraph range.)\n$openAiParagraph3 = $d.Paragraphs.Item($openAiIndex)\n$r3 = $openAiParagraph3.Range\n[void]$r3.MoveEnd(1, -1)\n[void]$r3.MoveEnd(1, -1 * $bookmarkSuffix.Length)\n$r3.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $r3)\n"}
